# Regenerate the handoff/handback report: the source UUID-named file
# "2fee83f0-a044-4832-a264-a08982a0b73b.md" was replaced by a freshly
# generated "de2f0e6c-f087-4fa4-9c6f-c804fd303f52.md", and the xliff
# artifact hashes / timestamps were refreshed accordingly.

$wb = $excel.ActiveWorkbook

$newId = "de2f0e6c-f087-4fa4-9c6f-c804fd303f52"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("G2").Value = "2016-08-18 00:53:55"

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newId.md"
}

# Column A was re-autosized by the report generator.
$wsOverview.Columns.Item(1).ColumnWidth = 38.48

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newId.md"
$wsZhCn.Range("G2").Value = "$newId.1fc4f996a9086de4e1176e20ec32355e52f818d3.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-18 00:53:50"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = "$newId.md"
}

$wsZhCn.Columns.Item(1).ColumnWidth = 38.48

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newId.md"
$wsDeDe.Range("G2").Value = "$newId.1fc4f996a9086de4e1176e20ec32355e52f818d3.de-de.xlf"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = "$newId.md"
}

$wsDeDe.Columns.Item(1).ColumnWidth = 38.48
